# Prefix reference-control ids ("Mxx") with "5G-" across the workbook.
#
# 1. In the "reference_controls" sheet, column A (ref_id) holds values like
#    "M1".."M144" (row 1 is the header "ref_id"). Each becomes "5G-M1".."5G-M144".
# 2. In the "requirements" sheet, column F (reference_controls) holds
#    comma-separated tokens like "1:M1,1:M2,1:M3" that list which reference
#    controls apply. Every "Mxx" token inside those strings needs the same
#    "5G-" prefix, e.g. "1:M1,1:M2" -> "1:5G-M1,1:5G-M2".

$wb = $excel.ActiveWorkbook

# --- reference_controls sheet: update ref_id column (A) ---
$refSheet = $wb.Worksheets.Item("reference_controls")
$refUsed = $refSheet.UsedRange
$refRows = $refUsed.Rows.Count

for ($r = 2; $r -le $refRows; $r++) {
    $cell = $refSheet.Cells.Item($r, 1)
    $text = $cell.Text
    if ($text -and $text.Length -gt 0) {
        $cell.Value = ($text -replace "M", "5G-M")
    }
}

# --- requirements sheet: update reference_controls column (F) ---
$reqSheet = $wb.Worksheets.Item("requirements")
$reqUsed = $reqSheet.UsedRange
$reqRows = $reqUsed.Rows.Count

for ($r = 2; $r -le $reqRows; $r++) {
    $cell = $reqSheet.Cells.Item($r, 6)
    $text = $cell.Text
    if ($text -and $text.Length -gt 0) {
        $cell.Value = ($text -replace "M", "5G-M")
    }
}
